$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.673.96'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.058.68'
$ws.Range("E3").Value = '  +3.42%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '385.10'
$ws.Range("E5").Value = '  +1.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.45'
$ws.Range("E6").Value = '  +1.24%  '
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -0.80%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.11'
$ws.Range("E10").Value = '  +2.02%  '
$ws.Range("E12").Value = '  +0.64%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.543.45'
$ws.Range("E13").Value = '  +3.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.71'
$ws.Range("E14").Value = '  +2.46%  '
$ws.Range("E15").Value = '  -0.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.065.33'
$ws.Range("E16").Value = '  +4.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.976'
$ws.Range("E17").Value = '  -1.70%  '
$ws.Range("E18").Value = '  -4.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '51.720.85'
$ws.Range("E19").Value = '  +1.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.15'
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.46'
$ws.Range("E21").Value = '  +0.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0965'
$ws.Range("E22").Value = '  +0.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.33'
$ws.Range("E23").Value = '  +0.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '268.93'
$ws.Range("E24").Value = '  +0.86%  '
$ws.Range("E25").Value = '  -1.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.44'
$ws.Range("E26").Value = '  +7.96%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '27.05'
$ws.Range("E27").Value = '  +4.62%  '
$ws.Range("E28").Value = '  +5.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.28'
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("E31").Value = '  -1.92%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.31'
$ws.Range("E32").Value = '  +0.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.50'
$ws.Range("E33").Value = '  +0.36%  '
$ws.Range("E34").Value = '  +0.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.48'
$ws.Range("E35").Value = '  -1.25%  '
$ws.Range("E36").Value = '  +2.84%  '
$ws.Range("E37").Value = '  -0.11%  '
$ws.Range("E38").Value = '  +4.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.292'
$ws.Range("E39").Value = '  +8.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.09'
$ws.Range("E40").Value = '  +3.88%  '
$ws.Range("E41").Value = '  +2.91%  '
$ws.Range("E42").Value = '  +3.20%  '
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("E44").Value = '  +1.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.73'
$ws.Range("E45").Value = '  +5.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.96'
$ws.Range("E46").Value = '  +2.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.52'
$ws.Range("E47").Value = '  +5.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.08'
$ws.Range("E48").Value = '  +3.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.046.88'
$ws.Range("E49").Value = '  +0.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.364.84'
$ws.Range("E50").Value = '  +3.48%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.206'
$ws.Range("E51").Value = '  +7.12%  '
